$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update title / header text (shared string rich-text runs) ---
$ws.Range("A8").Value = "Volume 32   Number  3"
$ws.Range("C9").Value = "Report Covering the Week  1/13/2025  Through  1/19/2025"

# --- Update weekly crime statistics table (rows 14-28) ---
# Row 14
$ws.Range("L14").Value = -50
$ws.Range("N14").Value = -75

# Row 15
$ws.Range("F14").Copy()
$ws.Range("D15").PasteSpecial(-4122)
$ws.Range("D15").Value = 1
$ws.Range("H14").Copy()
$ws.Range("E15").PasteSpecial(-4122)
$ws.Range("E15").Value = -100
$ws.Range("F15").Value = 2
$ws.Range("F14").Copy()
$ws.Range("G15").PasteSpecial(-4122)
$ws.Range("G15").Value = 1
$ws.Range("H14").Copy()
$ws.Range("H15").PasteSpecial(-4122)
$ws.Range("H15").Value = 100
$ws.Range("F14").Copy()
$ws.Range("J15").PasteSpecial(-4122)
$ws.Range("J15").Value = 1
$ws.Range("H14").Copy()
$ws.Range("K15").PasteSpecial(-4122)
$ws.Range("K15").Value = 0

# Row 16
$ws.Range("C16").Value = 6
$ws.Range("D16").Value = 5
$ws.Range("E16").Value = 20
$ws.Range("F16").Value = 17
$ws.Range("G16").Value = 22
$ws.Range("H16").Value = -22.727272727272
$ws.Range("I16").Value = 12
$ws.Range("J16").Value = 13
$ws.Range("K16").Value = -7.692307692307
$ws.Range("L16").Value = 33.333333333333
$ws.Range("M16").Value = 9.090909090909
$ws.Range("N16").Value = -67.567567567567

# Row 17
$ws.Range("C17").Value = 7
$ws.Range("D17").Value = 9
$ws.Range("E17").Value = -22.222222222222
$ws.Range("F17").Value = 39
$ws.Range("G17").Value = 26
$ws.Range("H17").Value = 50
$ws.Range("I17").Value = 25
$ws.Range("J17").Value = 21
$ws.Range("K17").Value = 19.047619047619
$ws.Range("L17").Value = 47.058823529411
$ws.Range("M17").Value = 316.666666666667
$ws.Range("N17").Value = -30.555555555555

# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D18").PasteSpecial(-4122)
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E18").PasteSpecial(-4122)
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 7
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 6
$ws.Range("K18").Value = 50
$ws.Range("L18").Value = -40
$ws.Range("M18").Value = 100
$ws.Range("N18").Value = -66.666666666666

# Row 19
$ws.Range("C19").Value = 5
$ws.Range("D19").Value = 9
$ws.Range("E19").Value = -44.444444444444
$ws.Range("G19").Value = 28
$ws.Range("H19").Value = -25
$ws.Range("I19").Value = 15
$ws.Range("J19").Value = 22
$ws.Range("K19").Value = -31.818181818181
$ws.Range("L19").Value = -34.782608695652
$ws.Range("M19").Value = 7.142857142857
$ws.Range("N19").Value = -59.459459459459

# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0"
$ws.Range("C14").Copy()
$ws.Range("D20").PasteSpecial(-4122)
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "***.*"
$ws.Range("C14").Copy()
$ws.Range("E20").PasteSpecial(-4122)
$ws.Range("F20").Value = 2
$ws.Range("G20").Value = 5
$ws.Range("H20").Value = -60

# Row 21
$ws.Range("C21").Value = 21
$ws.Range("D21").Value = 24
$ws.Range("E21").Value = -12.5
$ws.Range("F21").Value = 89
$ws.Range("G21").Value = 90
$ws.Range("H21").Value = -1.111111111111
$ws.Range("I21").Value = 60
$ws.Range("J21").Value = 64
$ws.Range("K21").Value = -6.25
$ws.Range("L21").Value = -13.043478260869
$ws.Range("M21").Value = 57.894736842105
$ws.Range("N21").Value = -61.538461538461

# Row 22
$ws.Range("F14").Copy()
$ws.Range("D22").PasteSpecial(-4122)
$ws.Range("D22").Value = 1
$ws.Range("H14").Copy()
$ws.Range("E22").PasteSpecial(-4122)
$ws.Range("E22").Value = -100
$ws.Range("F14").Copy()
$ws.Range("J22").PasteSpecial(-4122)
$ws.Range("J22").Value = 1
$ws.Range("H14").Copy()
$ws.Range("K22").PasteSpecial(-4122)
$ws.Range("K22").Value = -100

# Row 23
$ws.Range("C23").Value = 7
$ws.Range("D23").Value = 7
$ws.Range("E23").Value = 0
$ws.Range("F23").Value = 25
$ws.Range("G23").Value = 29
$ws.Range("H23").Value = -13.793103448275
$ws.Range("I23").Value = 13
$ws.Range("J23").Value = 20
$ws.Range("K23").Value = -35
$ws.Range("L23").Value = -40.90909090909
$ws.Range("M23").Value = 62.5

# Row 24
$ws.Range("C24").Value = 18
$ws.Range("D24").Value = 16
$ws.Range("E24").Value = 12.5
$ws.Range("G24").Value = 57
$ws.Range("H24").Value = 5.263157894736
$ws.Range("I24").Value = 37
$ws.Range("J24").Value = 33
$ws.Range("K24").Value = 12.121212121212
$ws.Range("L24").Value = 12.121212121212
$ws.Range("M24").Value = 2.777777777777

# Row 25
$ws.Range("C25").Value = 3
$ws.Range("F14").Copy()
$ws.Range("D25").PasteSpecial(-4122)
$ws.Range("D25").Value = 3
$ws.Range("H14").Copy()
$ws.Range("E25").PasteSpecial(-4122)
$ws.Range("E25").Value = 0
$ws.Range("G25").Value = 10
$ws.Range("H25").Value = -30
$ws.Range("I25").Value = 7
$ws.Range("J25").Value = 6
$ws.Range("K25").Value = 16.666666666666
$ws.Range("L25").Value = -22.222222222222

# Row 26
$ws.Range("C26").Value = 6
$ws.Range("D26").Value = 4
$ws.Range("E26").Value = 50
$ws.Range("F26").Value = 28
$ws.Range("G26").Value = 36
$ws.Range("H26").Value = -22.222222222222
$ws.Range("I26").Value = 21
$ws.Range("J26").Value = 20
$ws.Range("K26").Value = 5
$ws.Range("L26").Value = 0
$ws.Range("M26").Value = -32.258064516129

# Row 27
$ws.Range("F14").Copy()
$ws.Range("D27").PasteSpecial(-4122)
$ws.Range("D27").Value = 1
$ws.Range("H14").Copy()
$ws.Range("E27").PasteSpecial(-4122)
$ws.Range("E27").Value = -100
$ws.Range("F27").Value = 2
$ws.Range("F14").Copy()
$ws.Range("G27").PasteSpecial(-4122)
$ws.Range("G27").Value = 1
$ws.Range("H14").Copy()
$ws.Range("H27").PasteSpecial(-4122)
$ws.Range("H27").Value = 100
$ws.Range("F14").Copy()
$ws.Range("J27").PasteSpecial(-4122)
$ws.Range("J27").Value = 1
$ws.Range("H14").Copy()
$ws.Range("K27").PasteSpecial(-4122)
$ws.Range("K27").Value = 0

# Row 28
$ws.Range("F14").Copy()
$ws.Range("C28").PasteSpecial(-4122)
$ws.Range("C28").Value = 1
$ws.Range("F14").Copy()
$ws.Range("D28").PasteSpecial(-4122)
$ws.Range("D28").Value = 2
$ws.Range("H14").Copy()
$ws.Range("E28").PasteSpecial(-4122)
$ws.Range("E28").Value = -50
$ws.Range("F14").Copy()
$ws.Range("G28").PasteSpecial(-4122)
$ws.Range("G28").Value = 2
$ws.Range("H14").Copy()
$ws.Range("H28").PasteSpecial(-4122)
$ws.Range("H28").Value = 0
$ws.Range("F14").Copy()
$ws.Range("I28").PasteSpecial(-4122)
$ws.Range("I28").Value = 1
$ws.Range("F14").Copy()
$ws.Range("J28").PasteSpecial(-4122)
$ws.Range("J28").Value = 2
$ws.Range("H14").Copy()
$ws.Range("K28").PasteSpecial(-4122)
$ws.Range("K28").Value = -50
$ws.Range("L28").Value = -66.666666666666

$excel.CutCopyMode = 0

